$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.992.67'
$ws.Range("E2").Value = '  -4.04%  '
$ws.Range("D3").Value = '2.226.89'
$ws.Range("E3").Value = '  -6.34%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = "'296.30"
$ws.Range("E5").Value = '  -5.04%  '
$ws.Range("D6").Value = "'79.98"
$ws.Range("E6").Value = '  -8.45%  '
$ws.Range("D7").Value = "'0.507"
$ws.Range("E7").Value = '  -4.01%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = "'0.459"
$ws.Range("E9").Value = '  -6.39%  '
$ws.Range("D10").Value = "'0.0769"
$ws.Range("E10").Value = '  -6.38%  '
$ws.Range("D11").Value = "'27.78"
$ws.Range("E11").Value = '  -10.19%  '
$ws.Range("D12").Value = "'46.03"
$ws.Range("E12").Value = '  -13.08%  '
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").Value = '2.575.87'
$ws.Range("D15").Value = "'6.09"
$ws.Range("E15").Value = '  -7.84%  '
$ws.Range("D16").Value = "'14.03"
$ws.Range("E16").Value = '  -6.94%  '
$ws.Range("D17").Value = '2.240.55'
$ws.Range("E17").Value = '  -6.09%  '
$ws.Range("D18").Value = "'0.714"
$ws.Range("E18").Value = '  -5.56%  '
$ws.Range("D19").Value = '38.934.22'
$ws.Range("E19").Value = '  -3.84%  '
$ws.Range("D20").Value = '0.0₃0858'
$ws.Range("E20").Value = '  -5.81%  '
$ws.Range("D21").Value = "'5.72"
$ws.Range("E21").Value = '  -7.27%  '
$ws.Range("D22").Value = "'64.86"
$ws.Range("E22").Value = '  -5.90%  '
$ws.Range("D23").Value = "'9.76"
$ws.Range("E23").Value = '  -9.27%  '
$ws.Range("D24").Value = "'224.37"
$ws.Range("E24").Value = '  -4.90%  '
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("D26").Value = "'2.38"
$ws.Range("E26").Value = '  -9.55%  '
$ws.Range("E27").Value = '  -4.87%  '
$ws.Range("D28").Value = "'22.19"
$ws.Range("E28").Value = '  -5.61%  '
$ws.Range("D29").Value = "'2.13"
$ws.Range("E29").Value = '  -3.47%  '
$ws.Range("D30").Value = "'8.89"
$ws.Range("E30").Value = '  -5.10%  '
$ws.Range("D31").Value = "'149.28"
$ws.Range("E31").Value = '  -4.00%  '
$ws.Range("D32").Value = "'31.04"
$ws.Range("E32").Value = '  -8.08%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").Value = "'4.75"
$ws.Range("E34").Value = '  -8.87%  '
$ws.Range("E35").Value = '  -3.53%  '
$ws.Range("D36").Value = "'0.0682"
$ws.Range("E36").Value = '  -6.28%  '
$ws.Range("E37").Value = '  -4.29%  '
$ws.Range("D38").Value = "'2.65"
$ws.Range("E38").Value = '  -5.35%  '
$ws.Range("D39").Value = "'0.0952"
$ws.Range("E39").Value = '  -3.83%  '
$ws.Range("D40").Value = "'14.53"
$ws.Range("E40").Value = '  -8.47%  '
$ws.Range("D41").Value = "'1.60"
$ws.Range("E41").Value = '  -7.72%  '
$ws.Range("D42").Value = "'3.62"
$ws.Range("E42").Value = '  -5.29%  '
$ws.Range("D43").Value = '1.900.50'
$ws.Range("E43").Value = '  -2.97%  '
$ws.Range("E44").Value = '  -8.33%  '
$ws.Range("D45").Value = "'0.0253"
$ws.Range("E45").Value = '  -6.16%  '
$ws.Range("D46").Value = "'16.49"
$ws.Range("E46").Value = '  -6.14%  '
$ws.Range("E47").Value = '  -3.23%  '
$ws.Range("D48").Value = "'2.51"
$ws.Range("E48").Value = '  -10.45%  '
$ws.Range("D49").Value = '2.443.56'
$ws.Range("E49").Value = '  -6.20%  '
$ws.Range("D50").Value = "'87.38"
$ws.Range("E50").Value = '  -6.61%  '
$ws.Range("D51").Value = "'67.21"
$ws.Range("E51").Value = '  -7.62%  '
